$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update summary header counts (row 1) ---
$ws.Range("A1").Value = "Green Status: 14 projects"
$ws.Range("C1").Value = "Amber Status: 2 projects"

# --- Update existing data cells ---
$ws.Range("C2").Value = "Project 123: 14"
$ws.Range("E2").Value = "Project Twinkle Toes: 16"
$ws.Range("G2").Value = "c1: 17"

$ws.Range("C3").Value = "a3: 15"
$ws.Range("G3").Value = "c2: 18"

$ws.Range("G4").Value = "c3: 19"

$ws.Range("G5").Value = "o1: 20"

$ws.Range("G6").Value = "o2: 21"

$ws.Range("G7").Value = "a5: 22"

# --- Add new rows 13-15, cloning the A/C/E/G formatting from row 12 ---
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("C12").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("E12").Copy()
$ws.Range("E13").PasteSpecial(-4122)
$ws.Range("G12").Copy()
$ws.Range("G13").PasteSpecial(-4122)
$ws.Range("A13").Value = "a1: 11"

$ws.Range("A12").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("C12").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("E12").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("G12").Copy()
$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("A14").Value = "a2: 12"

$ws.Range("A12").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("C12").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("E12").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("G12").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("A15").Value = "a4: 13"

$excel.CutCopyMode = 0
